# "Final Import Form HR and running Power BI"
#
# The interview-tracking header row has three label cells (D1:F1) that were
# originally entered with a stray leading space (" InterviewResult",
# " HireDate", " InterviewerId"). Clean those up so the headers read
# correctly (this also causes Excel to re-pack the shared-string table,
# moving the Pass/Fail/Cancel result labels ahead of the corrected header
# labels - the same shared strings keep being used by the D2:D11 result
# cells, just addressed by their new index).
#
# Also move the active selection to F14, reflecting where the user's cursor
# ended up after finishing the import/review pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "InterviewResult"
$ws.Range("E1").Value = "HireDate"
$ws.Range("F1").Value = "InterviewerId"

$ws.Range("F14").Select()
